$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.810.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "'3.102.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'565.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'147.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'3.099.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("E10").Value = "  +3.93%  "
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("D12").Value = "'0.490"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.03%  "
$ws.Range("D14").Value = "'36.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("D15").Value = "'3.608.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "'64.734.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "'3.102.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "'7.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'499.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.01%  "
$ws.Range("D22").Value = "'15.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.16%  "
$ws.Range("D23").Value = "'0.709"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.85%  "
$ws.Range("D24").Value = "'7.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").Value = "'84.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'2.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.89%  "
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").Value = "'27.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "'2.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.87%  "
$ws.Range("D33").Value = "'1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "'6.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.68%  "
$ws.Range("D35").Value = "'6.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.39%  "
$ws.Range("D36").Value = "'55.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'458.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0417"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").Value = "'0.0839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("D40").Value = "'3.089.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("D41").Value = "'2.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.26%  "
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").Value = "'0.287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.63%  "
$ws.Range("D45").Value = "'2.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.94%  "
$ws.Range("D46").Value = "'28.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").Value = "'0.0₃0539"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("E50").Value = "  +6.60%  "
$ws.Range("D51").Value = "'117.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
